# Append a new data row (row 2) to the AfterShip returns-in-transit sheet,
# matching columns: tracking_number, carrier_slug, status_tag, order_id,
# last_checkpoint_id, last_checkpoint_time, last_checkpoint_location,
# updated_at, title.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "TEST_GDPR"
$ws.Range("B2").Value = "dbschenker-se"
$ws.Range("C2").Value = "Delivered"

# order_id is textual ("12345") in the source data, not numeric - force text
# entry with a leading apostrophe (quote-prefix) so Excel keeps it as a
# string instead of auto-converting to a number, then reset the cell style
# back to Normal so no stray quote-prefix formatting persists.
$ws.Range("D2").Value = "'12345"
$ws.Range("D2").Style = "Normal"

# last_checkpoint_id has no value for this shipment - leave it blank but
# still present as part of the row.
$ws.Range("E2").Font.Bold = $false

$ws.Range("F2").Value = "2026-02-08T16:46:34+04:30"
$ws.Range("G2").Value = ",Afghanistan"
$ws.Range("H2").Value = "2026-02-08T12:25:31+00:00"
$ws.Range("I2").Value = "SHIPMENT_TITLE"
